# "add chart setDataRange() example"
#
# On the "chart" worksheet:
#   - insert a new header row above the existing Chrome/IE/Firefox/... data
#     (data shifts from A1:B6 down to A2:B7) and label the two columns
#     "Column1" / "Column2"
#   - add a second, unrelated sample data block at A10:B16 (fruit / count)
#     with the same "Column1" / "Column2" headers
#   - turn both blocks into real Excel Tables (ListObjects) with distinct
#     styles, matching Table2 (TableStyleMedium9) and Table24
#     (TableStyleMedium10)
#   - widen column B a bit and leave the B11:B16 block selected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("chart")
$ws.Activate()

# Push the existing data down one row and add the Table2 header.
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Column1"
$ws.Range("B1").Value = "Column2"

# Second sample table: fruit counts at A10:B16.
$ws.Range("A10").Value = "Column1"
$ws.Range("B10").Value = "Column2"
$ws.Range("A11").Value = "Banna"
$ws.Range("B11").Value = 9
$ws.Range("A12").Value = "Apple"
$ws.Range("B12").Value = 11
$ws.Range("A13").Value = "Melon"
$ws.Range("B13").Value = 12
$ws.Range("A14").Value = "Cherry"
$ws.Range("B14").Value = 35
$ws.Range("A15").Value = "Mango"
$ws.Range("B15").Value = 22
$ws.Range("A16").Value = "Pineapple"
$ws.Range("B16").Value = 11

# Turn both ranges into Excel Tables.
$tbl1 = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:B7"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl2 = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A10:B16"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)

# Rename tbl2 before tbl1 so the auto-generated "Table2" default name (the
# second table Excel creates) never collides with the final name we want
# to give the first table.
$tbl2.Name = "Table24"
$tbl1.Name = "Table2"

$tbl1.TableStyle = "TableStyleMedium9"
$tbl2.TableStyle = "TableStyleMedium10"

# Cosmetic touch-ups matching the author's edit.
$ws.Columns.Item(2).ColumnWidth = 10.5
$ws.Range("B11:B16").Select()
